$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 80,
# pushing rows 80:113 down to 81:114 (dimension grows from T113 to T114).
$ws.Rows("80:80").Insert()

$ws.Range("A80").Value = 4
$ws.Range("B80").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C80").Value = "Los Lagos"
$ws.Range("D80").Value = 44466
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100101
$ws.Range("H80").Value = "Berries"
$ws.Range("I80").Value = 100101007
$ws.Range("J80").Value = "Kiwi"
$ws.Range("K80").Value = "Hayward"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 250
$ws.Range("N80").Value = 15000
$ws.Range("O80").Value = 15000
$ws.Range("P80").Value = 15000
$ws.Range("Q80").Value = "$/caja 15 kilos"
$ws.Range("R80").Value = "Provincia de Curicó"
$ws.Range("S80").Value = 1000
$ws.Range("T80").Value = 15
